$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "Magic mirror" (魔镜) row entirely - it was row 8.
# This shifts rows 9-11 up to become rows 8-10.
$ws.Rows(8).Delete()

# Fix up the style of the last row's effect cell (D10, "Blueprint" effect text)
# so it matches the normal wrap-less style used by the other D-column cells,
# instead of the stray style that row 11 used to carry.
$ws.Range("D9").Copy()
$ws.Range("D10").PasteSpecial(-4122)

# Update maxCount values (column C) per the card-balance changes:
# attribute / loot cards go from 4 to 3 max copies.
$ws.Range("C2").Value = 3   # 材料包 / Food pouch
$ws.Range("C3").Value = 3   # 钱袋 / Pouch
$ws.Range("C5").Value = 3   # 绷带 / Bandage
$ws.Range("C6").Value = 3   # 壶 / Pot

# Cursed variants now have an explicit maxCount of 1 (previously blank).
$ws.Range("C4").Value = 1   # 诅咒金币 / Cursed coin
$ws.Range("C7").Value = 1   # 诅咒之壶 / Cursed pot

# Restore the sheet view to the default top-left cell and update the
# active selection to C10.
[void]$ws.Range("C10").Select()
